$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'321.25"
$ws.Range("E2").Formula = "'7.76%"
$ws.Range("D3").Formula = "'49.10"
$ws.Range("E3").Formula = "'17.27%"
$ws.Range("D4").Formula = "'5.292"
$ws.Range("E4").Formula = "'5.60%"
$ws.Range("D5").Formula = "'0.08116"
$ws.Range("E5").Formula = "'7.70%"
$ws.Range("D6").Formula = "'4.612"
$ws.Range("E6").Formula = "'5.44%"
$ws.Range("D7").Formula = "'1.659"
$ws.Range("E7").Formula = "'3.37%"
$ws.Range("D8").Formula = "'1.206"
$ws.Range("E8").Formula = "'31.04%"
$ws.Range("D9").Formula = "'0.1321"
$ws.Range("E9").Formula = "'11.83%"
$ws.Range("D10").Formula = "'0.1951"
$ws.Range("E10").Formula = "'6.59%"
$ws.Range("D11").Formula = "'0.09548"
$ws.Range("E11").Formula = "'6.28%"
$ws.Range("D12").Formula = "'0.04515"
$ws.Range("E12").Formula = "'12.06%"
$ws.Range("D13").Formula = "'0.1049"
$ws.Range("E13").Formula = "'-0.01%"
$ws.Range("D14").Formula = "'0.001325"
$ws.Range("E14").Formula = "'3.42%"
$ws.Range("D15").Formula = "'0.005831"
$ws.Range("E15").Formula = "'-0.77%"
$ws.Range("E16").Formula = "'0.71%"
$ws.Range("E17").Formula = "'1.49%"
$ws.Range("D18").Formula = "'0.3389"
$ws.Range("E18").Formula = "'1.80%"
$ws.Range("D19").Formula = "'8.268"
$ws.Range("E19").Formula = "'-0.13%"
$ws.Range("D20").Formula = "'0.1410"
$ws.Range("E20").Formula = "'2.83%"
$ws.Range("D21").Formula = "'0.2922"
$ws.Range("E21").Formula = "'-9.27%"
$ws.Range("D22").Formula = "'0.04311"
$ws.Range("E22").Formula = "'5.33%"
$ws.Range("D23").Formula = "'0.001308"
$ws.Range("E23").Formula = "'3.32%"
$ws.Range("E24").Formula = "'9.03%"
$ws.Range("D25").Formula = "'0.0001355"
$ws.Range("E25").Formula = "'8.32%"
$ws.Range("D26").Formula = "'0.0003540"
$ws.Range("E26").Formula = "'-4.92%"
$ws.Range("D38").Formula = "'0.02724"
$ws.Range("E38").Formula = "'13.03%"
$ws.Range("D39").Formula = "'0.05596"
$ws.Range("E39").Formula = "'7.17%"
$ws.Range("D40").Formula = "'0.006322"
$ws.Range("E40").Formula = "'0.31%"
$ws.Range("D41").Formula = "'0.007703"
$ws.Range("D42").Formula = "'0.1443"
$ws.Range("E42").Formula = "'8.85%"
$ws.Range("D43").Formula = "'0.007698"
$ws.Range("E43").Formula = "'4.07%"
$ws.Range("E44").Formula = "'14.38%"
$ws.Range("E45").Formula = "'7.52%"
$ws.Range("D46").Formula = "'0.00006993"
$ws.Range("E46").Formula = "'6.12%"
$ws.Range("D47").Formula = "'0.00000000753"
$ws.Range("E47").Formula = "'0.31%"
$ws.Range("E48").Formula = "'35.00%"
$ws.Range("D49").Formula = "'0.004001"
$ws.Range("E49").Formula = "'-4.76%"
$ws.Range("D50").Formula = "'0.00002107"
$ws.Range("E50").Formula = "'0.31%"
$ws.Range("D51").Formula = "'0.0002007"
$ws.Range("E51").Formula = "'0.31%"
